$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 10112011
$ws.Range("K2").Value = "RuneStone"
$ws.Range("L2").Value = "Runestone"

$ws.Columns.Item(11).ColumnWidth = 9.285714285714286

$ws.Range("C8").Select()
